$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q0)
$ws.Range("B2").Value = -0.1679391225927079
$ws.Range("C2").Value = 0.6512668612981908
$ws.Range("D2").Value = 1.077333483753495
$ws.Range("E2").Value = 1.037946763448634
$ws.Range("F2").Value = 1.047290627843072
$ws.Range("G2").Value = 23

# Row 3 (Q1)
$ws.Range("B3").Value = 0.6108509255840642
$ws.Range("C3").Value = 0.8713162755896118
$ws.Range("D3").Value = 1.94510215081865
$ws.Range("E3").Value = 1.394669190460107
$ws.Range("F3").Value = 1.283284338193333
$ws.Range("G3").Value = 22

# Row 4 (Q2)
$ws.Range("B4").Value = 0.5135409635362268
$ws.Range("C4").Value = 1.149401313682296
$ws.Range("D4").Value = 3.430862853712828
$ws.Range("E4").Value = 1.852258851703192
$ws.Range("F4").Value = 1.8235941048126
$ws.Range("G4").Value = 21

# Row 5 (Q3)
$ws.Range("B5").Value = 0.6385223778103771
$ws.Range("C5").Value = 0.8102001327720327
$ws.Range("D5").Value = 1.083726471416818
$ws.Range("E5").Value = 1.041021840028737
$ws.Range("F5").Value = 0.8435611508437559
$ws.Range("G5").Value = 20

# Row 6 (Q4)
$ws.Range("B6").Value = 0.5184854322233537
$ws.Range("C6").Value = 0.6702212091134594
$ws.Range("D6").Value = 0.6945837070328604
$ws.Range("E6").Value = 0.8334168866976841
$ws.Range("F6").Value = 0.6703802697182577
$ws.Range("G6").Value = 19

# Row 7 (Q5)
$ws.Range("B7").Value = 0.3259858714718396
$ws.Range("C7").Value = 0.5159630503235952
$ws.Range("D7").Value = 0.3780268743366828
$ws.Range("E7").Value = 0.6148389011250693
$ws.Range("F7").Value = 0.5364195870263428
$ws.Range("G7").Value = 18

# Row 8 (Q6)
$ws.Range("B8").Value = 0.345936562473689
$ws.Range("C8").Value = 0.3960995140440492
$ws.Range("D8").Value = 0.2072861800771614
$ws.Range("E8").Value = 0.45528692060849
$ws.Range("F8").Value = 0.3051064642012103
$ws.Range("G8").Value = 17

# Row 9 (Q7)
$ws.Range("B9").Value = 0.3283187899062386
$ws.Range("C9").Value = 0.3609349132557007
$ws.Range("D9").Value = 0.1816422354358933
$ws.Range("E9").Value = 0.4261950673528418
$ws.Range("F9").Value = 0.2838354343252694
$ws.Range("G9").Value = 12

# Row 10 (Q8)
$ws.Range("B10").Value = 0.3306938847573825
$ws.Range("C10").Value = 0.3306938847573825
$ws.Range("D10").Value = 0.1666290365104273
$ws.Range("E10").Value = 0.4082022005212947
$ws.Range("F10").Value = 0.2584873103466553
$ws.Range("G10").Value = 7
